$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the Implementation-status tracker columns (P,Q) three columns to the
# right (to S,T) by inserting three blank columns at P.
$ws.Columns("P:R").Insert()

# Shift everything from row 28 down to row 30 by inserting two blank rows
# above the old row 28 (so its content lands on row 30).
$ws.Rows("28:29").Insert()

# Row 22 (FR-U-01 "화면 비례 길이 최소 단위 계산") requirement entry removed;
# only the category marker in D22 remains.
$ws.Range("H22").ClearContents()
$ws.Range("J22").ClearContents()
# The implementation-tracker cell for that row no longer tracks a status,
# so clear its color back to "No Fill".
$ws.Range("S22").Interior.ColorIndex = -4142

# Row 25 (FR-U-A-01) is now marked complete (green) instead of waiting (orange).
$ws.Range("S25").Interior.Color = 5287936

# Row 26 (FR-U-A-02): description text updated, and it now has a waiting
# (orange) status marker that it didn't have before.
$ws.Range("J26").Value = "게임 오브젝트 수직/수평 으로 입력한 거리만큼 일정/가속 속도로 이동"
$ws.Range("S26").Interior.Color = 49407

# Row 30 (formerly row 28, FR-U-T-01): description text updated to the new
# requirement about the notice/system message.
$ws.Range("J30").Value = "공지 메시지 = 시스템 메시지 출력"

# Cosmetic selection state to match the author's saved view.
$ws.Range("G18").Select() | Out-Null
